$d = $word.ActiveDocument

# The abstract paragraph currently reads (rendered text):
#   ... a NetHack public server (link). The data I discuss in this blog post is ...
# where "link" is a hyperlink to https://alt.org/nethack/. This edit drops the
# parenthetical hyperlink entirely, leaving:
#   ... a NetHack public server. The data I discuss in this blog post is ...

# Locate the "link" hyperlink robustly (by display text + target address)
# rather than assuming a fixed collection index.
$target = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $candidate = $d.Hyperlinks($i)
    if ($candidate.Range.Text -eq "link" -and $candidate.Address -eq "https://alt.org/nethack/") {
        $target = $candidate
        break
    }
}
if ($target -eq $null) {
    for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
        $candidate = $d.Hyperlinks($i)
        if ($candidate.Range.Text -eq "link") {
            $target = $candidate
            break
        }
    }
}

if ($target -ne $null) {
    # Strip the hyperlink field/formatting but keep its "link" text in place for
    # now, so every other run in the document is left completely untouched.
    $target.Delete()
}

# Collapse "public server (link)" down to "public server" - this removes the
# (now plain-text, if the hyperlink step above ran) "link" word along with its
# surrounding parentheses, joining the text before/after into one run.
$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute("public server (link)", $true, $false, $false, $false, $false, $true, 1, $false, "public server", 2)
Write-Output "Replaced 'public server (link)' -> 'public server': $found"
